$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows before the existing row 191, shifting the
# old rows 191-200 down to 194-203 (their values stay unchanged).
$ws.Range("A191:A193").EntireRow.Insert()

# Fill in the 3 newly inserted rows with the new observations.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are constant across this block of rows.

# Row 191
$ws.Range("A191").Value = 11
$ws.Range("B191").Value = "Vega Monumental Concepción"
$ws.Range("C191").Value = "Bíobío"
$ws.Range("D191").Value = 44516
$ws.Range("E191").Value = 8
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100101
$ws.Range("H191").Value = "Berries"
$ws.Range("I191").Value = 100112025
$ws.Range("J191").Value = "Frutilla"
$ws.Range("K191").Value = "Sin especificar"
$ws.Range("L191").Value = "Especial"
$ws.Range("M191").Value = 200
$ws.Range("N191").Value = 8000
$ws.Range("O191").Value = 8000
$ws.Range("P191").Value = 8000
$ws.Range("Q191").Value = '$/bandeja 7 kilos'
$ws.Range("R191").Value = "Provincia de Melipilla"
$ws.Range("S191").Value = 1143
$ws.Range("T191").Value = 7

# Row 192
$ws.Range("A192").Value = 11
$ws.Range("B192").Value = "Vega Monumental Concepción"
$ws.Range("C192").Value = "Bíobío"
$ws.Range("D192").Value = 44516
$ws.Range("E192").Value = 8
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100101
$ws.Range("H192").Value = "Berries"
$ws.Range("I192").Value = 100112025
$ws.Range("J192").Value = "Frutilla"
$ws.Range("K192").Value = "Sin especificar"
$ws.Range("L192").Value = "Primera"
$ws.Range("M192").Value = 1500
$ws.Range("N192").Value = 6500
$ws.Range("O192").Value = 6500
$ws.Range("P192").Value = 6500
$ws.Range("Q192").Value = '$/bandeja 7 kilos'
$ws.Range("R192").Value = "Provincia de Melipilla"
$ws.Range("S192").Value = 929
$ws.Range("T192").Value = 7

# Row 193
$ws.Range("A193").Value = 11
$ws.Range("B193").Value = "Vega Monumental Concepción"
$ws.Range("C193").Value = "Bíobío"
$ws.Range("D193").Value = 44516
$ws.Range("E193").Value = 8
$ws.Range("F193").Value = "Fruta"
$ws.Range("G193").Value = 100101
$ws.Range("H193").Value = "Berries"
$ws.Range("I193").Value = 100112025
$ws.Range("J193").Value = "Frutilla"
$ws.Range("K193").Value = "Sin especificar"
$ws.Range("L193").Value = "Segunda"
$ws.Range("M193").Value = 100
$ws.Range("N193").Value = 5000
$ws.Range("O193").Value = 5000
$ws.Range("P193").Value = 5000
$ws.Range("Q193").Value = '$/bandeja 7 kilos'
$ws.Range("R193").Value = "Provincia de Melipilla"
$ws.Range("S193").Value = 714
$ws.Range("T193").Value = 7
